$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "AC2" = "Bee"
    "AC4" = "Bee"
    "AC6" = "Bee"
    "AC7" = "Bee"
    "AC14" = "Bee"
    "AC15" = "Bee"
    "AC16" = "Bee"
    "AC18" = "Bee"
    "AC21" = "Bee"
    "AC22" = "Bee"
    "AC23" = "Butterfly"
    "AC24" = "Butterfly"
    "AC26" = "Bee"
    "AC27" = "Bee"
    "AC28" = "Bee"
    "AC29" = "Bee"
    "AC30" = "Bee"
    "AC31" = "Bee"
    "AC32" = "Bee"
    "AC33" = "Bee"
    "AC34" = "Bee"
    "AC35" = "Butterfly"
    "AC36" = "Bee"
    "AC37" = "Butterfly"
    "AC38" = "Bee"
    "AC40" = "Butterfly"
    "AC41" = "Bee"
    "AC45" = "Butterfly"
    "AC46" = "Butterfly"
    "AC48" = "Bee"
    "AC49" = "Bee"
    "AC51" = "Bee"
    "AC52" = "Bee"
    "AC53" = "Bee"
    "AC54" = "Bee"
    "AC56" = "Bee"
    "AC58" = "Bee"
    "AC59" = "Bee"
    "AC60" = "Bee"
    "AC61" = "Bee"
    "AC63" = "Bee"
    "AC64" = "Bee"
    "AC65" = "Bee"
    "AC67" = "Bee"
    "AC68" = "Bee"
    "AC69" = "Bee"
    "AC70" = "Bee"
    "AC74" = "Bee"
    "AC75" = "Bee"
    "AC78" = "Bee"
    "AC79" = "Butterfly"
    "AC80" = "Bee"
    "AC82" = "Bee"
    "AC83" = "Bee"
    "AC85" = "Bee"
    "AC86" = "Bee"
    "AC87" = "Bee"
    "AC88" = "Bee"
    "AC89" = "Bee"
    "AC92" = "Bee"
    "AC93" = "Bee"
    "AC94" = "Bee"
    "AC95" = "Bee"
    "AC96" = "Bee"
    "AC97" = "Bee"
    "AC98" = "Bee"
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
